{"js": "// Replace each two-digit-divided-by-one-digit answer string with its\n// updated value. Every answer string is unique within the document, so\n// searching for the exact literal text and replacing the single match is\n// unambiguous and leaves the surrounding run formatting (font/size) intact.\n\nconst pairs = [\n  [\"82\u00f76=13, 4\", \"94\u00f72=47, 0\"],\n  [\"44\u00f73=14, 2\", \"18\u00f72=9, 0\"],\n  [\"61\u00f72=30, 1\", \"56\u00f73=18, 2\"],\n  [\"95\u00f77=13, 4\", \"39\u00f75=7, 4\"],\n  [\"97\u00f73=32, 1\", \"76\u00f72=38, 0\"],\n  [\"86\u00f74=21, 2\", \"74\u00f72=37, 0\"],\n  [\"49\u00f72=24, 1\", \"18\u00f75=3, 3\"],\n  [\"78\u00f72=39, 0\", \"20\u00f73=6, 2\"],\n  [\"69\u00f79=7, 6\", \"20\u00f79=2, 2\"],\n  [\"46\u00f78=5, 6\", \"40\u00f78=5, 0\"],\n  [\"53\u00f75=10, 3\", \"67\u00f74=16, 3\"],\n  [\"98\u00f75=19, 3\", \"59\u00f78=7, 3\"],\n  [\"33\u00f74=8, 1\", \"60\u00f78=7, 4\"],\n  [\"22\u00f75=4, 2\", \"51\u00f75=10, 1\"],\n  [\"48\u00f72=24, 0\", \"41\u00f77=5, 6\"],\n  [\"99\u00f78=12, 3\", \"88\u00f79=9, 7\"],\n  [\"76\u00f73=25, 1\", \"47\u00f78=5, 7\"],\n  [\"35\u00f75=7, 0\", \"50\u00f73=16, 2\"],\n  [\"68\u00f78=8, 4\", \"17\u00f76=2, 5\"],\n  [\"12\u00f75=2, 2\", \"72\u00f74=18, 0\"],\n  [\"56\u00f74=14, 0\", \"67\u00f74=16, 3\"],\n  [\"20\u00f77=2, 6\", \"44\u00f72=22, 0\"],\n  [\"32\u00f73=10, 2\", \"53\u00f77=7, 4\"],\n  [\"67\u00f77=9, 4\", \"47\u00f74=11, 3\"],\n  [\"21\u00f78=2, 5\", \"13\u00f72=6, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit-divided-by-one-digit answer string with its\n# updated value. Every answer in the table is unique in the document, so a\n# plain literal Find/Replace (no wildcards) for each pair is unambiguous\n# and preserves each run's original formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"82\u00f76=13, 4\", \"94\u00f72=47, 0\"),\n    @(\"44\u00f73=14, 2\", \"18\u00f72=9, 0\"),\n    @(\"61\u00f72=30, 1\", \"56\u00f73=18, 2\"),\n    @(\"95\u00f77=13, 4\", \"39\u00f75=7, 4\"),\n    @(\"97\u00f73=32, 1\", \"76\u00f72=38, 0\"),\n    @(\"86\u00f74=21, 2\", \"74\u00f72=37, 0\"),\n    @(\"49\u00f72=24, 1\", \"18\u00f75=3, 3\"),\n    @(\"78\u00f72=39, 0\", \"20\u00f73=6, 2\"),\n    @(\"69\u00f79=7, 6\", \"20\u00f79=2, 2\"),\n    @(\"46\u00f78=5, 6\", \"40\u00f78=5, 0\"),\n    @(\"53\u00f75=10, 3\", \"67\u00f74=16, 3\"),\n    @(\"98\u00f75=19, 3\", \"59\u00f78=7, 3\"),\n    @(\"33\u00f74=8, 1\", \"60\u00f78=7, 4\"),\n    @(\"22\u00f75=4, 2\", \"51\u00f75=10, 1\"),\n    @(\"48\u00f72=24, 0\", \"41\u00f77=5, 6\"),\n    @(\"99\u00f78=12, 3\", \"88\u00f79=9, 7\"),\n    @(\"76\u00f73=25, 1\", \"47\u00f78=5, 7\"),\n    @(\"35\u00f75=7, 0\", \"50\u00f73=16, 2\"),\n    @(\"68\u00f78=8, 4\", \"17\u00f76=2, 5\"),\n    @(\"12\u00f75=2, 2\", \"72\u00f74=18, 0\"),\n    @(\"56\u00f74=14, 0\", \"67\u00f74=16, 3\"),\n    @(\"20\u00f77=2, 6\", \"44\u00f72=22, 0\"),\n    @(\"32\u00f73=10, 2\", \"53\u00f77=7, 4\"),\n    @(\"67\u00f77=9, 4\", \"47\u00f74=11, 3\"),\n    @(\"21\u00f78=2, 5\", \"13\u00f72=6, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
